$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $ref, $val)
    $c = $ws.Range($ref)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-TextValue $ws 'D2' '26.846.86'
Set-TextValue $ws 'E2' '  +0.45%  '
Set-TextValue $ws 'D3' '1.646.87'
Set-TextValue $ws 'E3' '  -0.10%  '
Set-TextValue $ws 'E4' '  +1.12%  '
Set-TextValue $ws 'D5' '217.62'
Set-TextValue $ws 'E5' '  +0.95%  '
Set-TextValue $ws 'E6' '  -0.12%  '
Set-TextValue $ws 'E7' '  +1.08%  '
Set-TextValue $ws 'D8' '0.251'
Set-TextValue $ws 'E8' '  -0.61%  '
Set-TextValue $ws 'E9' '  -0.34%  '
Set-TextValue $ws 'D10' '19.20'
Set-TextValue $ws 'E10' '  -0.72%  '
Set-TextValue $ws 'D11' '0.0843'
Set-TextValue $ws 'E11' '  -0.36%  '
Set-TextValue $ws 'D12' '1.871.67'
Set-TextValue $ws 'E12' '  -0.39%  '
Set-TextValue $ws 'D13' '1.652.53'
Set-TextValue $ws 'E13' '  -2.36%  '
Set-TextValue $ws 'D14' '4.18'
Set-TextValue $ws 'E14' '  -0.47%  '
Set-TextValue $ws 'D15' '0.527'
Set-TextValue $ws 'D16' '64.76'
Set-TextValue $ws 'E16' '  -2.30%  '
Set-TextValue $ws 'D17' '26.848.10'
Set-TextValue $ws 'E17' '  +0.14%  '
Set-TextValue $ws 'D18' '0.0₃0737'
Set-TextValue $ws 'E18' '  -2.49%  '
Set-TextValue $ws 'D19' '213.78'
Set-TextValue $ws 'E19' '  -3.44%  '
Set-TextValue $ws 'E20' '  +1.20%  '
Set-TextValue $ws 'D21' '4.35'
Set-TextValue $ws 'E21' '  -1.24%  '
Set-TextValue $ws 'E22' '  +11.49%  '
Set-TextValue $ws 'D23' '6.28'
Set-TextValue $ws 'E23' '  -1.86%  '
Set-TextValue $ws 'D24' '9.36'
Set-TextValue $ws 'E24' '  -1.98%  '
Set-TextValue $ws 'D25' '145.26'
Set-TextValue $ws 'E25' '  -2.00%  '
Set-TextValue $ws 'E26' '  +0.97%  '
Set-TextValue $ws 'E27' '  -2.54%  '
Set-TextValue $ws 'D28' '7.08'
Set-TextValue $ws 'D29' '15.69'
Set-TextValue $ws 'E29' '  -1.51%  '
Set-TextValue $ws 'D30' '0.0514'
Set-TextValue $ws 'E30' '  -1.59%  '
Set-TextValue $ws 'E31' '  +0.92%  '
Set-TextValue $ws 'D32' '3.31'
Set-TextValue $ws 'E32' '  -3.71%  '
Set-TextValue $ws 'E33' '  -2.73%  '
Set-TextValue $ws 'D34' '1.278.28'
Set-TextValue $ws 'E34' '  -1.71%  '
Set-TextValue $ws 'E35' '  -2.54%  '
Set-TextValue $ws 'D36' '2.45'
Set-TextValue $ws 'E36' '  +1.77%  '
Set-TextValue $ws 'D37' '0.0175'
Set-TextValue $ws 'E37' '  -4.09%  '
Set-TextValue $ws 'E38' '  +2.30%  '
Set-TextValue $ws 'D39' '0.825'
Set-TextValue $ws 'E39' '  +0.07%  '
Set-TextValue $ws 'E40' '  +1.13%  '
Set-TextValue $ws 'D41' '0.813'
Set-TextValue $ws 'E41' '  -0.28%  '
Set-TextValue $ws 'E42' '  -0.75%  '
Set-TextValue $ws 'E43' '  -1.46%  '
Set-TextValue $ws 'D44' '1.797.53'
Set-TextValue $ws 'E44' '  +0.40%  '
Set-TextValue $ws 'D45' '91.58'
Set-TextValue $ws 'E45' '  -2.29%  '
Set-TextValue $ws 'D46' '59.00'
Set-TextValue $ws 'E46' '  -3.52%  '
Set-TextValue $ws 'D47' '1.60'
Set-TextValue $ws 'E47' '  -0.96%  '
Set-TextValue $ws 'E48' '  -1.98%  '
Set-TextValue $ws 'E49' '  +0.27%  '
Set-TextValue $ws 'E50' '  -1.83%  '
Set-TextValue $ws 'D51' '0.0973'
Set-TextValue $ws 'E51' '  -0.42%  '
